# The source data set gained one new weekly price observation for
# "Zanahoria" (carrot) at Terminal Hortofrutícola Agro Chillán, recorded
# 2023-06-16 (serial 45093). New rows are inserted at the top of the
# date-ordered block (row 343), pushing every subsequent record down by
# one row (343->344, ..., 461->462) while rows 1-342 stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 343; this shifts rows
# 343:461 down to 344:462 (and grows the used range to R462), exactly
# like pressing "Insert Sheet Rows" in the Excel UI.
$ws.Rows.Item(343).Insert()

# Populate the newly inserted row 343 with the new observation. Columns
# A, B, C, E, F, G, H, I, N, O, Q, R keep the same values the row below
# (old row 343, now row 344) already has - only the date/volume/price
# columns (D, J, K, L, M, P) carry new data.
$ws.Range("A343").Value = 7
$ws.Range("B343").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C343").Value = "Ñuble"
$ws.Range("D343").Value = 45093
$ws.Range("E343").Value = 16
$ws.Range("F343").Value = 100114013
$ws.Range("G343").Value = "Zanahoria"
$ws.Range("H343").Value = "Sin especificar"
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 150
$ws.Range("K343").Value = 7000
$ws.Range("L343").Value = 7000
$ws.Range("M343").Value = 7000
$ws.Range("N343").Value = "`$/saco 20 kilos"
$ws.Range("O343").Value = "Región de Ñuble"
$ws.Range("P343").Value = 350
$ws.Range("Q343").Value = 20
$ws.Range("R343").Value = "Hortaliza"
